$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction of the two "unnamed" placeholder headers left over from the
# pandas multi-level column export: both should read "total", matching the
# other "total" sub-header in row 2.
$ws.Cells.Item(2, 2).Value = "total"
$ws.Cells.Item(2, 6).Value = "total"
